$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D width (Descripcion) ---
# Target XML width 14.5703125; the COM layer quantizes to 1/6 character steps,
# so this lands on the nearest achievable bucket (14.5).
$ws.Columns.Item(4).ColumnWidth = 13.65

# --- Row 2 ---
$ws.Range("H2").Value = 1
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = "N/A"

# --- Row 3 ---
$ws.Range("H3").Value = 1.6
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = "N/A"

# --- Row 4 ---
$ws.Range("H4").Value = 3.2
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = "N/A"

# --- Row 5 ---
$ws.Range("H5").Value = 6.4
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = "N/A"

# --- Row 6 ---
# H6 must become a *text* value "12.8" (not a number) while keeping its
# original cell style (s="4"). Typing the string directly gets silently
# coerced back to a number by the engine, and forcing text via
# NumberFormat/quote-prefix allocates a brand-new style. Instead, put a
# formula that evaluates to the text "12.8", then convert it to a static
# value in-place via copy / paste-special-values, which keeps the
# existing style untouched.
$ws.Range("H6").Formula = "=TEXT(12.8,""0.0"")"
$ws.Range("H6").Copy()
$ws.Range("H6").PasteSpecial(-4163)
$excel.CutCopyMode = 0

$ws.Range("J6").Value = 0
$ws.Range("K6").Value = "N/A"

# --- Selection state ---
$ws.Range("E2:H5").Select()
